# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The "dimension" mapping columns balance / ratios / rama-descripcion / pyg
# are recurated as plain measures (matching the other measure columns),
# while "ano" stays the only real dimension column. For each of the four
# affected columns (I, AN, BF, BH) this:
#   - demotes the iaest-dimension:* label (row 2) to iaest-measure:*
#   - flips the "dim" marker (row 3) to "medida"
#   - flips the "skos:Concept" type (row 4) to "xsd:int"
#   - removes the mapping-*.xlsx cell (row 5) entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("I", "AN", "BF", "BH")

foreach ($col in $cols) {
    # Row 2: iaest-dimension:<name> -> iaest-measure:<name>
    $label = $ws.Range($col + "2").Text
    $ws.Range($col + "2").Value = $label -replace "^iaest-dimension:", "iaest-measure:"

    # Row 3: dim -> medida
    $ws.Range($col + "3").Value = "medida"

    # Row 4: skos:Concept -> xsd:int
    $ws.Range($col + "4").Value = "xsd:int"

    # Row 5: drop the mapping-*.xlsx cell completely (content + style)
    $ws.Range($col + "5").Clear()
}
